# Refresh the forecast run: the generation script was re-run a week later,
# so every week's row slides forward by one week and all the forecast
# numbers (MyForecast / Amazon Mean / P70 / P80 / P90) are refreshed with
# the new model output. The Summary sheet's derived statistics are updated
# to match.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "2025-02-02"
$ws1.Range("D2").Value = 523
$ws1.Range("E2").Value = 402
$ws1.Range("F2").Value = 464
$ws1.Range("G2").Value = 517
$ws1.Range("H2").Value = 597

$ws1.Range("B3").NumberFormat = "@"
$ws1.Range("B3").Value = "2025-02-09"
$ws1.Range("D3").Value = 446
$ws1.Range("E3").Value = 343
$ws1.Range("F3").Value = 399
$ws1.Range("G3").Value = 448
$ws1.Range("H3").Value = 522

$ws1.Range("B4").NumberFormat = "@"
$ws1.Range("B4").Value = "2025-02-16"
$ws1.Range("D4").Value = 462
$ws1.Range("E4").Value = 355
$ws1.Range("F4").Value = 414
$ws1.Range("G4").Value = 466
$ws1.Range("H4").Value = 545

$ws1.Range("B5").NumberFormat = "@"
$ws1.Range("B5").Value = "2025-02-23"
$ws1.Range("D5").Value = 469
$ws1.Range("E5").Value = 361
$ws1.Range("F5").Value = 423
$ws1.Range("G5").Value = 479
$ws1.Range("H5").Value = 565

$ws1.Range("B6").NumberFormat = "@"
$ws1.Range("B6").Value = "2025-03-02"
$ws1.Range("D6").Value = 457
$ws1.Range("E6").Value = 357
$ws1.Range("F6").Value = 420
$ws1.Range("G6").Value = 480
$ws1.Range("H6").Value = 571

$ws1.Range("B7").NumberFormat = "@"
$ws1.Range("B7").Value = "2025-03-09"
$ws1.Range("D7").Value = 394
$ws1.Range("E7").Value = 362
$ws1.Range("F7").Value = 426
$ws1.Range("G7").Value = 485
$ws1.Range("H7").Value = 576

$ws1.Range("B8").NumberFormat = "@"
$ws1.Range("B8").Value = "2025-03-16"
$ws1.Range("D8").Value = 373
$ws1.Range("E8").Value = 354
$ws1.Range("F8").Value = 419
$ws1.Range("G8").Value = 482
$ws1.Range("H8").Value = 579

$ws1.Range("B9").NumberFormat = "@"
$ws1.Range("B9").Value = "2025-03-23"
$ws1.Range("D9").Value = 369
$ws1.Range("E9").Value = 368
$ws1.Range("F9").Value = 436
$ws1.Range("G9").Value = 502
$ws1.Range("H9").Value = 603

$ws1.Range("B10").NumberFormat = "@"
$ws1.Range("B10").Value = "2025-03-30"
$ws1.Range("D10").Value = 370
$ws1.Range("E10").Value = 353
$ws1.Range("F10").Value = 418
$ws1.Range("G10").Value = 482
$ws1.Range("H10").Value = 579

$ws1.Range("B11").NumberFormat = "@"
$ws1.Range("B11").Value = "2025-04-06"
$ws1.Range("D11").Value = 352
$ws1.Range("E11").Value = 346
$ws1.Range("F11").Value = 412
$ws1.Range("G11").Value = 476
$ws1.Range("H11").Value = 574

$ws1.Range("B12").NumberFormat = "@"
$ws1.Range("B12").Value = "2025-04-13"
$ws1.Range("D12").Value = 341
$ws1.Range("E12").Value = 343
$ws1.Range("F12").Value = 410
$ws1.Range("G12").Value = 478
$ws1.Range("H12").Value = 583

$ws1.Range("B13").NumberFormat = "@"
$ws1.Range("B13").Value = "2025-04-20"
$ws1.Range("D13").Value = 339
$ws1.Range("E13").Value = 339
$ws1.Range("F13").Value = 408
$ws1.Range("G13").Value = 481
$ws1.Range("H13").Value = 594

$ws1.Range("B14").NumberFormat = "@"
$ws1.Range("B14").Value = "2025-04-27"
$ws1.Range("D14").Value = 369
$ws1.Range("E14").Value = 333
$ws1.Range("F14").Value = 399
$ws1.Range("G14").Value = 467
$ws1.Range("H14").Value = 573

$ws1.Range("B15").NumberFormat = "@"
$ws1.Range("B15").Value = "2025-05-04"
$ws1.Range("D15").Value = 395
$ws1.Range("E15").Value = 324
$ws1.Range("F15").Value = 388
$ws1.Range("G15").Value = 455
$ws1.Range("H15").Value = 559

$ws1.Range("B16").NumberFormat = "@"
$ws1.Range("B16").Value = "2025-05-11"
$ws1.Range("D16").Value = 396
$ws1.Range("E16").Value = 321
$ws1.Range("F16").Value = 386
$ws1.Range("G16").Value = 455
$ws1.Range("H16").Value = 563

$ws1.Range("B17").NumberFormat = "@"
$ws1.Range("B17").Value = "2025-05-18"
$ws1.Range("D17").Value = 384
$ws1.Range("E17").Value = 314
$ws1.Range("F17").Value = 378
$ws1.Range("G17").Value = 447
$ws1.Range("H17").Value = 556

$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "2023-01-01 to 2025-01-26"

$ws2.Range("B4").NumberFormat = "@"
$ws2.Range("B4").Value = "1004"

$ws2.Range("B5").NumberFormat = "@"
$ws2.Range("B5").Value = "395"

$ws2.Range("B6").NumberFormat = "@"
$ws2.Range("B6").Value = "361"

$ws2.Range("B7").NumberFormat = "@"
$ws2.Range("B7").Value = "239"

$ws2.Range("B8").NumberFormat = "@"
$ws2.Range("B8").Value = "41872 units"

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "6438"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "3492"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "1899"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "523"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "339"

$ws2.Range("B15").NumberFormat = "@"
$ws2.Range("B15").Value = "2025-04-20"
